$wb = $excel.ActiveWorkbook

# Helper: write a literal text value into a cell (shared string), bypassing Excel's
# automatic "TRUE"/"FALSE" -> boolean literal coercion, and leaving the cell with the
# plain default style afterwards (no quote-prefix / text-format residue).
function Set-TextCell($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1. Insert two new validation-list sheets, "is_targeted list" and
#    "is_technical_replicate list", right after "analyte_class list" and before
#    "library_concentration_unit list" (so the final tab order becomes:
#    ... analyte_class list, is_targeted list, is_technical_replicate list,
#    library_concentration_unit list, library_final_yield_unit list, library_layout list)
# ---------------------------------------------------------------------------
$afterAnalyte = $wb.Worksheets.Item("analyte_class list")

$isTargetedSheet = $wb.Worksheets.Add($null, $afterAnalyte)
$isTargetedSheet.Name = "is_targeted list"

$isTechRepSheet = $wb.Worksheets.Add($null, $isTargetedSheet)
$isTechRepSheet.Name = "is_technical_replicate list"

# Populate both new list sheets with the two allowed boolean-as-text values.
Set-TextCell $isTargetedSheet.Cells.Item(1, 1) "TRUE"
Set-TextCell $isTargetedSheet.Cells.Item(2, 1) "FALSE"

Set-TextCell $isTechRepSheet.Cells.Item(1, 1) "TRUE"
Set-TextCell $isTechRepSheet.Cells.Item(2, 1) "FALSE"

# ---------------------------------------------------------------------------
# 2. Re-point the "is_targeted" (column N) and "is_technical_replicate" (column S)
#    data validations on the "Export as TSV" sheet away from the old inline
#    "TRUE,FALSE" literal-list formula and at the new list sheets instead, updating
#    the validation error text to match the "Value must come from list" convention
#    used by the other list-backed columns.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Export as TSV")

$isTargetedCol = $ws.Range("N2:N1048576")
$isTargetedCol.Validation.Formula1 = "'is_targeted list'!`$A`$1:`$A`$2"
$isTargetedCol.Validation.ErrorTitle = "Value must come from list"
$isTargetedCol.Validation.ErrorMessage = "Value must be one of: TRUE / FALSE."

$isTechRepCol = $ws.Range("S2:S1048576")
$isTechRepCol.Validation.Formula1 = "'is_technical_replicate list'!`$A`$1:`$A`$2"
$isTechRepCol.Validation.ErrorTitle = "Value must come from list"
$isTechRepCol.Validation.ErrorMessage = "Value must be one of: TRUE / FALSE."
